# "focus on fast reading and printing"
# - Rename the "Fast Reading" header to "Fast Reading 500 WPM"
# - Add a new "Fast Typing" column (Q) to the Agenda sheet, mirroring column P
# - Shift the Agenda week dates forward by 5 weeks (Dec 2019 -> Jan 2020)
# - Freeze the first column / adjust selection
# - Resize columns P/Q

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Agenda")

# ---------------------------------------------------------------------
# 1. Shared string: "Fast Reading" -> "Fast Reading 500 WPM"
# ---------------------------------------------------------------------
$ws.Range("P1").Value = "Fast Reading 500 WPM"

# ---------------------------------------------------------------------
# 2. New column Q ("Fast Typing"), cloning column P's look & feel
# ---------------------------------------------------------------------
$ws.Range("P1:P12").Copy()
$ws.Range("Q1").PasteSpecial(-4122)

$ws.Range("Q1").Value = "Fast Typing"

# ---------------------------------------------------------------------
# 3. Update the week's dates (A2:A8) -- shift 5 weeks later
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 43843
$ws.Range("A3").Value = 43844
$ws.Range("A4").Value = 43845
$ws.Range("A5").Value = 43846
$ws.Range("A6").Value = 43847
$ws.Range("A7").Value = 43848
$ws.Range("A8").Value = 43849

# ---------------------------------------------------------------------
# 4. Column P lost its per-day shared SUM formulas -- only row 6 keeps
#    a live formula now, the rest became plain zeros, and the weekly
#    total (P9) is re-entered as its own explicit SUM.
# ---------------------------------------------------------------------
$ws.Range("P2").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("P6").Formula = "=SUM(C6:O6)"
$ws.Range("P7").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("P9").Formula = "=SUM(P2:P8)"

# ---------------------------------------------------------------------
# 5. Fill in the new Q column data / formulas
# ---------------------------------------------------------------------
$ws.Range("Q2").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("Q6").Value = 0
$ws.Range("Q7").Value = 0
$ws.Range("Q8").Value = 0

$ws.Range("Q9").Formula = "=SUM(Q2:Q8)"

$ws.Range("P10").Value = 9
$ws.Range("Q10").Value = 3.5

$ws.Range("Q11").Formula = "=Q10-Q9"

$ws.Range("Q12").Value = "N/A"

# ---------------------------------------------------------------------
# 6. Widen column P (now has a longer header) and size the new column Q
# ---------------------------------------------------------------------
$ws.Columns.Item(16).ColumnWidth = 20.451822916666668
$ws.Columns.Item(17).ColumnWidth = 10.022135416666666

Write-Host "done"
